$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the weekly points totals (column B), and tag each row's
# working column D with the same "0.0" number format already used on H.
$ws.Range("B2").Value = 2245.5
$ws.Range("D2").NumberFormat = "0.0"

$ws.Range("B3").Value = 2201.5
$ws.Range("D3").NumberFormat = "0.0"

$ws.Range("B4").Value = 2079.3000000000002
$ws.Range("D4").NumberFormat = "0.0"

# Rows 5 & 6 swap athletes (Jeremiah Gaddy now ranks above Todd Vinsant)
$ws.Range("A5").Value = "Jeremiah Gaddy"
$ws.Range("B5").Value = 1645.7
$ws.Range("D5").NumberFormat = "0.0"

$ws.Range("A6").Value = "Todd Vinsant"
$ws.Range("B6").Value = 1590.7
$ws.Range("D6").NumberFormat = "0.0"

$ws.Range("B7").Value = 1390.1
$ws.Range("D7").NumberFormat = "0.0"

$ws.Range("B8").Value = 1117.3
$ws.Range("D8").NumberFormat = "0.0"

$ws.Range("D9").NumberFormat = "0.0"

$ws.Range("B10").Value = 811.1
$ws.Range("D10").NumberFormat = "0.0"

$ws.Range("D11").NumberFormat = "0.0"

# Move the active selection to F12
$ws.Range("F12").Select()
